# Update the header labels in row 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Replace the numeric RLS_option / n_clusters flag (column E) with the
# textual adaptive-filter name "RLS" for every data row, and refresh the
# recomputed RMSE / NDEI / MAE metrics (columns F, G, H) with their
# updated values.

$data = @(
    @{ Row = 2;  F = 0.8284200847349987; G = 3.693637885881631;  H = 0.655919892663091  },
    @{ Row = 3;  F = 0.2179103806544387; G = 0.971586822366329;  H = 0.1695656338448241 },
    @{ Row = 4;  F = 0.216632022992024;  G = 0.965887068846818;  H = 0.1681594209027977 },
    @{ Row = 5;  F = 0.1892557507211808; G = 0.8438257640848219; H = 0.1312410552030424 },
    @{ Row = 6;  F = 0.181328705168839;  G = 0.8084818168353956; H = 0.127594695458468  },
    @{ Row = 7;  F = 0.1763298313757094; G = 0.786193572055634;  H = 0.1251791102276751 },
    @{ Row = 8;  F = 0.1682753659500402; G = 0.750281503776541;  H = 0.1212318459407656 },
    @{ Row = 9;  F = 0.1569494668392412; G = 0.6997832471333522; H = 0.1154952920977721 },
    @{ Row = 10; F = 0.1435270275349208; G = 0.6399372447863136; H = 0.1083065630062123 },
    @{ Row = 11; F = 0.1305558347890536; G = 0.5821031943642572; H = 0.1006629752162537 },
    @{ Row = 12; F = 0.1201565341138518; G = 0.5357363188280526; H = 0.09357753230076804 },
    @{ Row = 13; F = 0.1120166337044638; G = 0.4994433255829205; H = 0.08886493754502671 },
    @{ Row = 14; F = 0.1038296379905924; G = 0.4629403507063769; H = 0.08355876720772551 },
    @{ Row = 15; F = 0.09341882198238849; G = 0.4165221322934979; H = 0.07543156409866535 },
    @{ Row = 16; F = 0.08531923223119302; G = 0.3804088702946845; H = 0.0690510800455076 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 5).Value = "RLS"
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
}
